$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rows 2-34: price/volume updates only ---
$ws.Range("D2").Value = "'70.959.10"
$ws.Range("E2").Value = "  +2.14%  "
$ws.Range("D3").Value = "'3.584.20"
$ws.Range("E3").Value = "  +1.07%  "
$ws.Range("D4").Value = "'1.00"
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").Value = "'606.16"
$ws.Range("E5").Value = "  +4.24%  "
$ws.Range("D6").Value = "'175.29"
$ws.Range("E6").Value = "  +1.26%  "
$ws.Range("D7").Value = "'3.577.90"
$ws.Range("E7").Value = "  +1.02%  "
$ws.Range("D8").Value = "'0.619"
$ws.Range("E8").Value = "  +0.89%  "
$ws.Range("E9").Value = "  -0.05%  "
$ws.Range("D10").Value = "'0.198"
$ws.Range("E10").Value = "  +3.85%  "
$ws.Range("E11").Value = "  +10.68%  "
$ws.Range("D12").Value = "'0.592"
$ws.Range("E12").Value = "  +1.11%  "
$ws.Range("D13").Value = "'47.21"
$ws.Range("E13").Value = "  -0.78%  "
$ws.Range("D14").Value = "'0.0000280"
$ws.Range("E14").Value = "  +1.14%  "
$ws.Range("D15").Value = "'4.164.20"
$ws.Range("E15").Value = "  +1.24%  "
$ws.Range("D16").Value = "'8.46"
$ws.Range("E16").Value = "  -1.28%  "
$ws.Range("D17").Value = "'621.95"
$ws.Range("E17").Value = "  -1.47%  "
$ws.Range("D18").Value = "'3.584.22"
$ws.Range("E18").Value = "  +1.30%  "
$ws.Range("D19").Value = "'71.012.24"
$ws.Range("E19").Value = "  +2.21%  "
$ws.Range("E20").Value = "  -2.56%  "
$ws.Range("D21").Value = "'17.53"
$ws.Range("E21").Value = "  +0.59%  "
$ws.Range("D22").Value = "'0.892"
$ws.Range("E22").Value = "  -0.26%  "
$ws.Range("D23").Value = "'9.42"
$ws.Range("E23").Value = "  -16.13%  "
$ws.Range("D24").Value = "'16.24"
$ws.Range("E24").Value = "  +1.41%  "
$ws.Range("D25").Value = "'97.97"
$ws.Range("E25").Value = "  +0.33%  "
$ws.Range("D26").Value = "'3.81"
$ws.Range("E26").Value = "  -0.03%  "
$ws.Range("E27").Value = "  +0.08%  "
$ws.Range("D28").Value = "'2.66"
$ws.Range("E28").Value = "  +0.63%  "
$ws.Range("D29").Value = "'9.40"
$ws.Range("E29").Value = "  +0.42%  "
$ws.Range("D30").Value = "'33.69"
$ws.Range("E30").Value = "  +2.03%  "
$ws.Range("D31").Value = "'8.55"
$ws.Range("E31").Value = "  -0.54%  "
$ws.Range("D32").Value = "'3.08"
$ws.Range("E32").Value = "  -2.54%  "
$ws.Range("D33").Value = "'7.17"
$ws.Range("E33").Value = "  +2.12%  "
$ws.Range("E34").Value = "  -2.77%  "

# --- Rows 35-51: new coin (Bittensor) inserted, rest shift down, Cronos dropped ---
$ws.Range("B35").Value = "Bittensor"
$ws.Range("C35").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D35").Value = "'648.18"
$ws.Range("E35").Value = "  +1.39%  "
$ws.Range("B36").Value = "dogwifhat"
$ws.Range("C36").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D36").Value = "'3.82"
$ws.Range("E36").Value = "  +8.06%  "
$ws.Range("B37").Value = "Hedera"
$ws.Range("C37").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D37").Value = "'0.103"
$ws.Range("E37").Value = "  -0.51%  "
$ws.Range("B38").Value = "Cosmos"
$ws.Range("C38").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D38").Value = "'10.92"
$ws.Range("E38").Value = "  +0.92%  "
$ws.Range("B39").Value = "VeChain"
$ws.Range("C39").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D39").Value = "'0.0489"
$ws.Range("E39").Value = "  +6.92%  "
$ws.Range("B40").Value = "OKB"
$ws.Range("C40").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D40").Value = "'57.60"
$ws.Range("E40").Value = "  +0.30%  "
$ws.Range("B41").Value = "FirstDigitalUSD"
$ws.Range("C41").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D41").Value = "'1.00"
$ws.Range("E41").Value = "  +0.05%  "
$ws.Range("B42").Value = "Kaspa"
$ws.Range("C42").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D42").Value = "'0.142"
$ws.Range("E42").Value = "  +3.68%  "
$ws.Range("B43").Value = "Maker"
$ws.Range("C43").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D43").Value = "'3.416.45"
$ws.Range("E43").Value = "  +0.47%  "
$ws.Range("B44").Value = "TheGraph"
$ws.Range("C44").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D44").Value = "'0.326"
$ws.Range("E44").Value = "  -1.12%  "
$ws.Range("B45").Value = "PEPE"
$ws.Range("C45").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D45").Value = "'0.0₃0721"
$ws.Range("E45").Value = "  +2.17%  "
$ws.Range("B46").Value = "ThetaToken"
$ws.Range("C46").Value = "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
$ws.Range("D46").Value = "'3.00"
$ws.Range("E46").Value = "  +8.85%  "
$ws.Range("B47").Value = "InjectiveProtocol"
$ws.Range("C47").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D47").Value = "'33.27"
$ws.Range("E47").Value = "  +1.22%  "
$ws.Range("B48").Value = "Fetch.AI"
$ws.Range("C48").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D48").Value = "'2.70"
$ws.Range("E48").Value = "  +4.88%  "
$ws.Range("B49").Value = "Stellar"
$ws.Range("C49").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D49").Value = "'0.131"
$ws.Range("E49").Value = "  +1.09%  "
$ws.Range("B50").Value = "Monero"
$ws.Range("C50").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D50").Value = "'133.04"
$ws.Range("E50").Value = "  -0.02%  "
$ws.Range("B51").Value = "USDe"
$ws.Range("C51").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D51").Value = "'1.00"
$ws.Range("E51").Value = "  -0.03%  "
